$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "51.237.73"
$ws.Range("E2").Value2 = "  -1.81%  "

# Row 3
$ws.Range("D3").Value2 = "2.917.39"
$ws.Range("E3").Value2 = "  -2.43%  "

# Row 4
$ws.Range("E4").Value2 = "  -0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "373.87"
$ws.Range("E5").Value2 = "  +5.52%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "103.78"
$ws.Range("E6").Value2 = "  -3.92%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.543"
$ws.Range("E7").Value2 = "  -3.47%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.999"
$ws.Range("E8").Value2 = "  -0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.591"
$ws.Range("E9").Value2 = "  -5.79%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "37.15"
$ws.Range("E10").Value2 = "  -3.19%  "

# Row 11
$ws.Range("E11").Value2 = "  -0.39%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.0838"
$ws.Range("E12").Value2 = "  -2.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "18.36"
$ws.Range("E13").Value2 = "  -5.29%  "

# Row 14
$ws.Range("D14").Value2 = "3.377.13"
$ws.Range("E14").Value2 = "  -2.43%  "

# Row 15
$ws.Range("E15").Value2 = "  -4.45%  "

# Row 16
$ws.Range("D16").Value2 = "2.916.18"
$ws.Range("E16").Value2 = "  -2.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.932"
$ws.Range("E17").Value2 = "  -9.08%  "

# Row 18
$ws.Range("D18").Value2 = "51.175.13"
$ws.Range("E18").Value2 = "  -2.00%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "3.32"
$ws.Range("E19").Value2 = "  -5.41%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "7.28"
$ws.Range("E20").Value2 = "  -3.80%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "13.02"
$ws.Range("E21").Value2 = "  -4.81%  "

# Row 22
$ws.Range("D22").Value2 = "0.0₃0945"
$ws.Range("E22").Value2 = "  -3.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "68.47"
$ws.Range("E23").Value2 = "  -1.54%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "260.64"
$ws.Range("E24").Value2 = "  -1.42%  "

# Row 25
$ws.Range("E25").Value2 = "  -1.81%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.171"
$ws.Range("E26").Value2 = "  -6.07%  "

# Row 27
$ws.Range("E27").Value2 = "  -4.02%  "

# Row 28
$ws.Range("E28").Value2 = "  -0.05%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "25.81"
$ws.Range("E29").Value2 = "  -4.02%  "

# Row 30
$ws.Range("E30").Value2 = "  -5.30%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "6.62"
$ws.Range("E31").Value2 = "  +4.45%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.102"
$ws.Range("E32").Value2 = "  -5.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "9.92"
$ws.Range("E33").Value2 = "  -4.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "34.60"
$ws.Range("E35").Value2 = "  -5.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "50.98"
$ws.Range("E36").Value2 = "  +0.22%  "

# Row 37
$ws.Range("E37").Value2 = "  +0.49%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.0425"
$ws.Range("E38").Value2 = "  -4.57%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "3.00"
$ws.Range("E39").Value2 = "  -6.24%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "17.08"
$ws.Range("E40").Value2 = "  -4.67%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "2.58"
$ws.Range("E41").Value2 = "  -4.56%  "

# Row 42
$ws.Range("E42").Value2 = "  -6.49%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.113"
$ws.Range("E43").Value2 = "  -3.44%  "

# Row 44
$ws.Range("B44").Value2 = "EnergySwap"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "21.97"
$ws.Range("E44").Value2 = "  -4.18%  "

# Row 45
$ws.Range("B45").Value2 = "Monero"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "119.40"
$ws.Range("E45").Value2 = "  -2.58%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "2.09"
$ws.Range("E46").Value2 = "  -2.78%  "

# Row 47
$ws.Range("D47").Value2 = "2.026.15"
$ws.Range("E47").Value2 = "  -4.72%  "

# Row 48
$ws.Range("E48").Value2 = "  -3.88%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "3.19"
$ws.Range("E49").Value2 = "  -6.05%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.251"
$ws.Range("E50").Value2 = "  +0.70%  "

# Row 51
$ws.Range("D51").Value2 = "3.214.46"
$ws.Range("E51").Value2 = "  -2.19%  "
